$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) values as plain text strings (quote-prefixed)
# to preserve the inlineStr/text representation of numeric-looking values.
$ws.Range("D2").Value = "'275.58"
$ws.Range("D3").Value = "'23.25"
$ws.Range("D4").Value = "'6.480"
$ws.Range("D5").Value = "'0.06293"
$ws.Range("D7").Value = "'6.683"
$ws.Range("D8").Value = "'1.395"
$ws.Range("D9").Value = "'0.8346"
$ws.Range("D10").Value = "'0.01384"
$ws.Range("D11").Value = "'0.1637"
$ws.Range("D12").Value = "'0.08319"
$ws.Range("D14").Value = "'0.03115"
$ws.Range("D15").Value = "'0.09308"
$ws.Range("D16").Value = "'3.872"
$ws.Range("D17").Value = "'0.001637"
$ws.Range("D18").Value = "'0.04785"
$ws.Range("D19").Value = "'0.006422"
$ws.Range("D20").Value = "'0.005693"
$ws.Range("D40").Value = "'0.04724"
$ws.Range("D41").Value = "'0.007049"
$ws.Range("D42").Value = "'0.1162"
$ws.Range("D43").Value = "'0.003700"
$ws.Range("D45").Value = "'0.00006239"
$ws.Range("D48").Value = "'0.7964"
$ws.Range("D50").Value = "'0.001988"
